$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column G: fill in the "Notes" results for rows 3-14 -------------------
$ws.Range("G3").Value  = "Reason: No matching Account Number on the lookup table. Account Number: TESTING"
$ws.Range("G4").Value  = "Reason: No matching Account Number on the lookup table. Account Number: 991046"
$ws.Range("G5").Value  = "Reason: No matching Account Number on the lookup table. Account Number: 991047"
$ws.Range("G6").Value  = "Reason: No matching Account Number on the lookup table. Account Number: 991048"
$ws.Range("G7").Value  = "Reason: No matching Account Number on the lookup table. Account Number: 991040"
$ws.Range("G8").Value  = "Reason: No matching Account Number on the lookup table. Account Number: 991045"
$ws.Range("G9").Value  = "Reason: No matching Account Number on the lookup table. Account Number: 991041"
$ws.Range("G10").Value = "Reason: No matching Account Number on the lookup table. Account Number: 991044"
$ws.Range("G11").Value = "Reason: No matching Account Number on the lookup table. Account Number: 991002"
$ws.Range("G12").Value = "Reason: No matching Account Number on the lookup table. Account Number: 533704"
$ws.Range("G13").Value = "Success, CCR"
$ws.Range("G14").Value = "Success, Assign To"

# --- Column A: new SR# for the last regression row --------------------------
$ws.Range("A14").Value = "SR0003015462"

# --- View state: zoom to 70% and move the selection to E14 ------------------
$ws.Select()
$excel.ActiveWindow.Zoom = 70
$ws.Range("E14").Select()
